# Commit: "remove spaces before first name initials"
#
# The "Initial" column (column B) of the AuthorList sheet contains many
# values that were stored with one or more leading spaces (e.g. " V.",
# " I.J.") while others have none (e.g. "C.", "H.").  This edit strips any
# leading spaces from every value in column B, leaving the rest of the
# string (including any trailing padding spaces) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2

    if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Length -gt 0 -and $val.Substring(0, 1) -eq " ") {
        $cell.Value2 = $val.TrimStart(" ")
    }
}

# Best-effort: restore the author's on-screen scroll position / selection
# (view-state only, does not affect any cell content).
$win = $excel.ActiveWindow
[void]$ws.Range("A76").Select()
$win.ScrollRow = 76
$win.ScrollColumn = 1
[void]$ws.Range("B113").Select()

Write-Host "Leading spaces removed from column B (Initial)."
